$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.172.79'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '1.862.14'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''0.7089'
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").Value = '''241.21'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").Value = '''0.9999'
$ws.Range("D8").Value = '''0.3101'
$ws.Range("E8").Value = '  -0.54%  '
$ws.Range("D9").Value = '''0.07652'
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").Value = '''24.71'
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").Value = '''0.08351'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("D12").Value = '1.866.30'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").Value = '''5.189'
$ws.Range("E13").Value = '  -2.03%  '
$ws.Range("D14").Value = '''0.7082'
$ws.Range("E14").Value = '  -2.78%  '
$ws.Range("D15").Value = '''91.19'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '29.172.39'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '''5.920'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '''243.01'
$ws.Range("E18").Value = '  -1.74%  '
$ws.Range("D19").Value = '''0.000007806'
$ws.Range("E19").Value = '  -0.83%  '
$ws.Range("D20").Value = '2.115.34'
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").Value = '''13.09'
$ws.Range("E21").Value = '  -1.87%  '
$ws.Range("D22").Value = '''0.9997'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '''7.872'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").Value = '''0.1581'
$ws.Range("E25").Value = '  -1.41%  '
$ws.Range("D26").Value = '''163.28'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").Value = '''8.960'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("D28").Value = '''18.43'
$ws.Range("E28").Value = '  +0.47%  '
$ws.Range("D29").Value = '''1.330'
$ws.Range("E29").Value = '  -2.36%  '
$ws.Range("E30").Value = '  -0.30%  '
$ws.Range("D31").Value = '''4.404'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '''4.237'
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").Value = '''0.05142'
$ws.Range("E33").Value = '  -2.74%  '
$ws.Range("D34").Value = '''0.7939'
$ws.Range("E34").Value = '  +9.06%  '
$ws.Range("E35").Value = '  -1.22%  '
$ws.Range("D36").Value = '''1.164'
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("D37").Value = '''2.683'
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").Value = '''0.01845'
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("E39").Value = '  -1.23%  '
$ws.Range("D40").Value = '1.166.74'
$ws.Range("E40").Value = '  -5.55%  '
$ws.Range("D41").Value = '''6.200'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''0.8894'
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("D43").Value = '''72.92'
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("D44").Value = '''0.9995'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").Value = '''102.03'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("D46").Value = '2.008.48'
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").Value = '''0.5198'
$ws.Range("E47").Value = '  -1.67%  '
$ws.Range("D48").Value = '''1.774'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '''0.00000000119'
$ws.Range("E50").Value = '  -0.09%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '''0.4275'
$ws.Range("E51").Value = '  -1.11%  '
